# Add a "time_taken" metadata column (F) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1: copy formatting from the existing E1 header (bold,
# bordered, centered style) so F1 matches the rest of the header row, then
# set its text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("F1").Value = "time_taken"

# Per-row time_taken timestamps for rows 2..103 (data rows), stored as plain
# text, in row order.
$timeTaken = @(
    "2021-10-05 10:50:24.589953",
    "2021-10-05 10:50:24.589964",
    "2021-10-05 10:50:24.589967",
    "2021-10-05 10:50:24.589970",
    "2021-10-05 10:50:24.589973",
    "2021-10-05 10:50:24.589975",
    "2021-10-05 10:50:24.589978",
    "2021-10-05 10:50:24.589980",
    "2021-10-05 10:50:24.589983",
    "2021-10-05 10:50:24.589985",
    "2021-10-05 10:50:24.589988",
    "2021-10-05 10:50:24.589990",
    "2021-10-05 10:50:24.589993",
    "2021-10-05 10:50:24.589995",
    "2021-10-05 10:50:24.589998",
    "2021-10-05 10:50:24.590000",
    "2021-10-05 10:50:24.590003",
    "2021-10-05 10:50:24.590006",
    "2021-10-05 10:50:24.590009",
    "2021-10-05 10:50:24.590011",
    "2021-10-05 10:50:24.590013",
    "2021-10-05 10:50:24.590016",
    "2021-10-05 10:50:24.590018",
    "2021-10-05 10:50:24.590021",
    "2021-10-05 10:50:24.590024",
    "2021-10-05 10:50:24.590026",
    "2021-10-05 10:50:24.590029",
    "2021-10-05 10:50:24.590031",
    "2021-10-05 10:50:24.590034",
    "2021-10-05 10:50:24.590036",
    "2021-10-05 10:50:24.590039",
    "2021-10-05 10:50:24.590041",
    "2021-10-05 10:50:24.590044",
    "2021-10-05 10:50:24.590046",
    "2021-10-05 10:50:24.590049",
    "2021-10-05 10:50:24.590052",
    "2021-10-05 10:50:24.590054",
    "2021-10-05 10:50:24.590057",
    "2021-10-05 10:50:24.590059",
    "2021-10-05 10:50:24.590062",
    "2021-10-05 10:50:24.590065",
    "2021-10-05 10:50:24.590067",
    "2021-10-05 10:50:24.590070",
    "2021-10-05 10:50:24.590072",
    "2021-10-05 10:50:24.590075",
    "2021-10-05 10:50:24.590077",
    "2021-10-05 10:50:24.590080",
    "2021-10-05 10:50:24.590082",
    "2021-10-05 10:50:24.590084",
    "2021-10-05 10:50:24.590087",
    "2021-10-05 10:50:24.590089",
    "2021-10-05 10:50:24.590092",
    "2021-10-05 10:50:24.590095",
    "2021-10-05 10:50:24.590097",
    "2021-10-05 10:50:24.590100",
    "2021-10-05 10:50:24.590102",
    "2021-10-05 10:50:24.590104",
    "2021-10-05 10:50:24.590107",
    "2021-10-05 10:50:24.590109",
    "2021-10-05 10:50:24.590112",
    "2021-10-05 10:50:24.590114",
    "2021-10-05 10:50:24.590117",
    "2021-10-05 10:50:24.590119",
    "2021-10-05 10:50:24.590122",
    "2021-10-05 10:50:24.590125",
    "2021-10-05 10:50:24.590128",
    "2021-10-05 10:50:24.590130",
    "2021-10-05 10:50:24.590133",
    "2021-10-05 10:50:24.590135",
    "2021-10-05 10:50:24.590138",
    "2021-10-05 10:50:24.590140",
    "2021-10-05 10:50:24.590143",
    "2021-10-05 10:50:24.590146",
    "2021-10-05 10:50:24.590148",
    "2021-10-05 10:50:24.590151",
    "2021-10-05 10:50:24.590153",
    "2021-10-05 10:50:24.590157",
    "2021-10-05 10:50:24.590160",
    "2021-10-05 10:50:24.590163",
    "2021-10-05 10:50:24.590165",
    "2021-10-05 10:50:24.590168",
    "2021-10-05 10:50:24.590170",
    "2021-10-05 10:50:24.590173",
    "2021-10-05 10:50:24.590175",
    "2021-10-05 10:50:24.590178",
    "2021-10-05 10:50:24.590180",
    "2021-10-05 10:50:24.590183",
    "2021-10-05 10:50:24.590185",
    "2021-10-05 10:50:24.590188",
    "2021-10-05 10:50:24.590190",
    "2021-10-05 10:50:24.590193",
    "2021-10-05 10:50:24.590195",
    "2021-10-05 10:50:24.590199",
    "2021-10-05 10:50:24.590202",
    "2021-10-05 10:50:24.590204",
    "2021-10-05 10:50:24.590207",
    "2021-10-05 10:50:24.590209",
    "2021-10-05 10:50:24.590212",
    "2021-10-05 10:50:24.590214",
    "2021-10-05 10:50:24.590217",
    "2021-10-05 10:50:24.590219",
    "2021-10-05 10:50:24.590222"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
